$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(87).Insert()

$ws.Cells.Item(87, 1).Value = 2
$ws.Cells.Item(87, 2).Value = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(87, 3).Value = "Coquimbo"
$ws.Cells.Item(87, 4).Value = 45106
$ws.Cells.Item(87, 5).Value = 4
$ws.Cells.Item(87, 6).Value = 100112026
$ws.Cells.Item(87, 7).Value = "Haba"
$ws.Cells.Item(87, 8).Value = "Sin especificar"
$ws.Cells.Item(87, 9).Value = "Primera"
$ws.Cells.Item(87, 10).Value = 1100
$ws.Cells.Item(87, 11).Value = 11000
$ws.Cells.Item(87, 12).Value = 12000
$ws.Cells.Item(87, 13).Value = 11500
$ws.Cells.Item(87, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(87, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(87, 16).Value = 460
$ws.Cells.Item(87, 17).Value = 25
$ws.Cells.Item(87, 18).Value = "Hortaliza"
